# Generate Report for Handoff
# - Flip status from "In Translation" to "Ready for handoff" on all three
#   sheets (Overview summary columns + per-locale detail sheets).
# - Bump the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps to the new handoff-generation time.
# - Columns that hold the (now longer) status text widen to fit.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status columns (E, F) + HO generate date (G)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-21 07:04:26"

# --- zh-cn detail sheet: Status column (C) + Latest Handoff Datetime (H)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-21 07:04:22"

# --- de-de detail sheet: Status column (C) + Latest Handoff Datetime (H)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-21 07:04:26"

# --- Column widths widen slightly to accommodate the new, longer status text
$wsOverview.Columns.Item(5).ColumnWidth = 16.34
$wsOverview.Columns.Item(6).ColumnWidth = 16.34
$wsZhCn.Columns.Item(3).ColumnWidth = 16.34
$wsDeDe.Columns.Item(3).ColumnWidth = 16.34
